$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the existing "brand name" cell style (Arial font, style index 2)
# from A3 onto the new rows' A cells before writing their values, so we
# reuse the existing style record instead of minting a new one.
$ws.Range("A3").Copy()
$ws.Range("A8:A10").PasteSpecial(-4122)

# Row 8 - new "Test" record
$ws.Range("A8").Value = "Test"
$ws.Range("B8").Value = 0
$ws.Range("C8").Value = 0
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = 0
$ws.Range("F8").Value = 0

# Row 9 - new "Test" record (D9 keeps its pre-existing underline style)
$ws.Range("A9").Value = "Test"
$ws.Range("B9").Value = 1
$ws.Range("C9").Value = 0
$ws.Range("D9").Value = 1
$ws.Range("E9").Value = 0
$ws.Range("F9").Value = 1

# Row 10 - new "Test" record
$ws.Range("A10").Value = "Test"
$ws.Range("B10").Value = 0
$ws.Range("C10").Value = 1
$ws.Range("D10").Value = 0
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 0

# Update the active selection to match the edited workbook
$ws.Range("E11").Select()
